# Add header row ("21th_day_test and 60th_day_test into Estrus" field names)
# to the "基本資料" sheet, shifting all existing data rows down by one, and
# fix up the two J-column messages that referenced the generic "{field}"
# placeholder so they now mention the concrete field name "Chinese_name".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基本資料")

# Insert a brand-new row at the top; this pushes all existing rows (and
# their values/number formats/styles) down by one without altering them.
$ws.Rows.Item(1).Insert()

# Populate the new header row with the English field names.
$headers = @("Breed", "ID", "confusing_note", "Birthday", "Sire", "Dam", "reg_id", "Chinese_name", "Gender", "註釋")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# The two rows that previously complained about a blank Chinese_name value
# used a generic "{field}" placeholder in column J; replace it with the
# actual field name now that the header row exists. After the row insert
# above, these data rows moved from 10/26 to 11/27.
$oldMessage = "['{field} 不能為空值', '不允許有相近耳號']"
$newMessage = "['Chinese_name 不能為空值', '不允許有相近耳號']"

foreach ($r in @(11, 27)) {
    $cell = $ws.Cells.Item($r, 10)
    if ($cell.Value() -eq $oldMessage) {
        $cell.Value = $newMessage
    }
}
